$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Terms fixed for customers" - the term labels in column B were renamed
$ws.Range("B1").Value = "terms1"
$ws.Range("B2").Value = "terms2"

# C2 picked up a stray one-off font; drop back to the regular font used
# elsewhere on the sheet (matches C1's formatting) instead of the special
# font that only this cell used.
$ws.Range("C2").Font.Name = "Arial"

# Active selection moved to B2
$ws.Range("B2").Select()
